$d = $word.ActiveDocument

# 1. "Ativacao" credits line: update activation date
$d.Content.Find.Execute('Ativação: 01/01/1996', $true, $false, $false, $false, $false,
                         $true, 1, $false, 'Ativação: 01/01/2022', 2) | Out-Null

# 2. "Objetivos" paragraph rewrite
$d.Content.Find.Execute('Fornecer aos estudantes uma visão abrangente e interdisciplinar dos materiais compostos por fases caracterizadas por distintos tipos de materiais (metais, cerâmicas e polímeros) para obter propriedades únicas. Apresentar os fundamentos teóricos da mecânica de estruturas reforçadas com fibras, tecidos e partículas. Apresentar os diferentes tipos de materiais compósitos, inclusive sobre os nanocompósitos e compósitos funcionais, que representam o avanço mais recente na área de Ciência e Engenharia de Materiais.', $true, $false, $false, $false, $false,
                         $true, 1, $false, 'Fornecer aos estudantes uma visão abrangente e interdisciplinar sobre materiais compósitos, além de mostrar as especificidades de cada matriz, sendo ela metálica, cerâmica ou polimérica. Ademais, deseja-se apresentar os fundamentos teóricos da mecânica de estruturas reforçadas e a partir de atividades práticas demostrar métodos de caracterização de materiais compósitos e como prepara-los.', 2) | Out-Null

# 3. Add two more docentes responsaveis after "519033 - Carlos Yujiro Shigue"
$profRange = $d.Paragraphs(8).Range
$profRange.Text = "519033 - Carlos Yujiro Shigue`v1033242 - Fábio Herbst Florenzano`v1922320 - Sebastiao Ribeiro"

# 4. "Programa resumido" paragraph rewrite
$d.Content.Find.Execute('Materiais compósitos: tipos, propriedades, processamento e aplicações. Nanocompósitos e compósitos funcionais.', $true, $false, $false, $false, $false,
                         $true, 1, $false, '1.Introduçâo 2. Conceitos básicos sobre materiais compósitos, suas matrizes e seus processo de fabricação 3. Tipos de reforços 4. Compósitos nanoestruturados, naturais e híbridos 5. Mecânica da estrutura reforçada 6. Atividade prática', 2) | Out-Null

# 5. "Programa" paragraph rewrite
$d.Content.Find.Execute('Conteúdo teórico:1. Conceitos básicos sobre materiais compósitos: compósitos de matriz metálica (CMM), compósitos de matriz cerâmicos (CMC) e compósitos de matriz polimérica (CMP) e nanocompósitos.2. Fibras, tecidos e reforços particulados.3. Mecânica de estruturas reforçadas.4. Compósitos de matriz metálica: características e processos de fabricação.5. Compósitos de matriz cerâmica: características e processos de fabricação.6. Compósitos de matriz polimérica: matrizes termoplásticas e termorrígidas, características físicas e químicas e processos de fabricação.7. Compósitos nanoestruturados.8.Compósitos funcionais.Conteúdo prático:1. Caracterização e análise de compósitos de matriz metálica.2. Preparação e caracterização de compósito de matriz polimérica.3. Visita a empresa produtora de compósitos.', $true, $false, $false, $false, $false,
                         $true, 1, $false, '1. Conceitos básicos sobre materiais compósitos: compósitos de matriz metálica (CMM), compósitos de matriz cerâmicos (CMC) e compósitos de matriz polimérica (CMP) e nanocompósitos. 2. Tipos de Reforços: Reforços particulados, fibras curtas, fibras longas, mantas, tecidos e preformas. 3. Conceitos de Interface4. Compósitos de matriz metálica: características e processos de fabricação. 5. Compósitos de matriz cerâmica: características e processos de fabricação. 6. Compósitos de matriz polimérica: matrizes termoplásticas e termorrígidas, características físicas e químicas e processos de fabricação. 7. Compósitos nanoestruturados. 8. Compósitos Naturais. 9. Compósitos Híbridos 10. Mecânica de estruturas reforçadas. Conteúdo prático: 1. Caracterização e análise de compósitos de matriz metálica. 2. Preparação e caracterização de compósitos de matriz polimérica.(Sugestão: Considerar substituir essa parte prática pela realização do PBL descrito no item 3) 3. Visita a empresa produtora de compósitos e aulas especiais e/ou palestras com professores/pesquisadores convidados', 2) | Out-Null

# 6. "Metodo" value rewrite (Avaliacao section)
$d.Content.Find.Execute('A avaliação será feita por meio de provas escritas.', $true, $false, $false, $false, $false,
                         $true, 1, $false, 'De acordo com a atual ementa da disciplina propõe-se o uso de uma nova metodologia de ensino com o intuito de abordar o conteúdo de forma mais prática e contextualizada para que o aluno consiga relacionar os conhecimentos teóricos vistos em sala de aula com as outras disciplinas do curso. Assim, avaliação do aluno será feita através de uma prova escrita e por uma apresentação final com base nas atividades práticas desenvolvidas.', 2) | Out-Null

# 7. "Criterio" value rewrite
$d.Content.Find.Execute('A Nota final (NF) será calculada da seguinte maneira:NF = (P1 + 2*P2)/3', $true, $false, $false, $false, $false,
                         $true, 1, $false, 'A nota final será calculada como descrita a seguir: NF= (0,4*Avaliação escrita + 0,6 *Apresentação final)', 2) | Out-Null

# 8. "Norma de recuperacao" value rewrite
$d.Content.Find.Execute('A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2', $true, $false, $false, $false, $false,
                         $true, 1, $false, 'Devido a cunho prático da disciplina não haverá recuperação.', 2) | Out-Null

# 9. "Bibliografia" paragraph rewrite
$d.Content.Find.Execute('1. MALLICK, P.K. Composites Engineering Handbook. New York: Marcel Dekker, 1997.2. MATTHEWS, F.L. & RAWLINGS, R.D. Composite Materials: Engineering and Science. London: Chapman & Hall, 1994.3. OBRAZTSOV, I.F. Mechanics of Composites. Moscow: MIR Publishers, 1982.4. JONES R. Mechanics of Composite Materials. New York: McGraw-Hill, 1975.5. UPADHYAYA, G.S. Sintered Metal-Ceramic Composites. Elsevier, 1984.6. HARPER, C. A. Handbook of Plastics, Elastomers and Composites. New York: McGraw-Hill, 1992.7. GOLDSTEIN, A.N. Handbook of Nanophase Materials. CRC Press, 1997.8. DRESSELHAUS, M.S. Graphite Fibers and Filaments. New York: Springer-Verlag, 1988.', $true, $false, $false, $false, $false,
                         $true, 1, $false, '1. REZENDE, M. C.; COSTA, M. L.; BOTELHO, E. C. Compósitos estruturais: tecnologia e prática. São Paulo: Artliber, 2011. 396p. 2 MALLICK, P.K. Composites Engineering Handbook. New York: Marcel Dekker, 1997. 3. MATTHEWS, F.L. & RAWLINGS, R.D. Composite Materials: Engineering and Science. London: Chapman & Hall, 1994. 4. OBRAZTSOV, I.F. Mechanics of Composites. Moscow: MIR Publishers, 1982. 5. JONES R. Mechanics of Composite Materials. New York: McGraw-Hill, 1975. 6. UPADHYAYA, G.S. Sintered Metal-Ceramic Composites. Elsevier, 1984. 7. HARPER, C. A. Handbook of Plastics, Elastomers and Composites. New York: McGraw-Hill, 1992. 8. GOLDSTEIN, A.N. Handbook of Nanophase Materials. CRC Press, 1997. 9. DRESSELHAUS, M.S. Graphite Fibers and Filaments. New York: Springer-Verlag, 1988.', 2) | Out-Null

